# Update cryptocurrency price/volume data as of Tue Sep 19 18:57:13 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.187.34"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.649.20"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D13").Value = "1.644.86"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "27.180.11"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").Value = "  +3.14%  "
$ws.Range("E22").Value = "  +5.26%  "
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("E26").Value = "  +3.15%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0508"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("D35").Value = "1.268.98"
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("E38").Value = "  +3.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.851"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("E43").Value = "  +6.23%  "
$ws.Range("D44").Value = "1.790.44"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0975"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "
